# Update "想去人数" (F column) values on sheets "展览" and "全部类型"
$wb = $excel.ActiveWorkbook

$updates = @{
    "展览"     = @{ "F3" = 271; "F4" = 81; "F5" = 273; "F6" = 270; "F7" = 99;
                     "F13" = 2428; "F14" = 41; "F16" = 13; "F18" = 531; "F19" = 571;
                     "F20" = 170; "F22" = 50; "F24" = 2004; "F25" = 4130; "F27" = 65;
                     "F28" = 1204; "F29" = 234; "F30" = 2111; "F32" = 472; "F35" = 293;
                     "F36" = 429; "F38" = 714; "F41" = 425 }
    "全部类型" = @{ "F3" = 271; "F4" = 81; "F5" = 273; "F6" = 270; "F7" = 99;
                     "F13" = 2428; "F14" = 41; "F17" = 13; "F19" = 531; "F20" = 571;
                     "F21" = 170; "F23" = 50; "F25" = 2004; "F26" = 4130; "F28" = 65;
                     "F29" = 1204; "F30" = 234; "F31" = 2111; "F33" = 472; "F36" = 293;
                     "F37" = 429; "F39" = 714; "F42" = 425 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellAddr in $cellUpdates.Keys) {
        $ws.Range($cellAddr).Value = $cellUpdates[$cellAddr]
    }
}

$wb.Save()
